$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1423.634542915053
$ws.Cells.Item(3, 3).Value = 1571.040212937525
$ws.Cells.Item(4, 3).Value = 1559.093043380422
$ws.Cells.Item(5, 3).Value = 1449.703934039648
$ws.Cells.Item(6, 3).Value = 1446.906377072232
$ws.Cells.Item(7, 3).Value = 1419.431552814099
$ws.Cells.Item(8, 3).Value = 1336.606502926552
$ws.Cells.Item(9, 3).Value = 1382.273267210077
$ws.Cells.Item(10, 3).Value = 1401.322065436084
$ws.Cells.Item(11, 3).Value = 1354.847459798858
$ws.Cells.Item(12, 3).Value = 1348.388858080751
$ws.Cells.Item(13, 3).Value = 1363.687312056716
$ws.Cells.Item(14, 3).Value = 1353.389188105791
$ws.Cells.Item(15, 3).Value = 1346.483021144484
$ws.Cells.Item(16, 3).Value = 1318.724169067753
$ws.Cells.Item(17, 3).Value = 1299.128641493841
$ws.Cells.Item(18, 3).Value = 1272.019449547579
$ws.Cells.Item(19, 3).Value = 1277.82174897606
$ws.Cells.Item(20, 3).Value = 1274.371438275579
$ws.Cells.Item(21, 3).Value = 1264.819362376879
$ws.Cells.Item(22, 3).Value = 1249.658422274535
$ws.Cells.Item(23, 3).Value = 1233.41148765066
$ws.Cells.Item(24, 3).Value = 1217.793156547686
$ws.Cells.Item(25, 3).Value = 1207.725347116329
$ws.Cells.Item(26, 3).Value = 1197.28210749781
$ws.Cells.Item(27, 3).Value = 1187.323709887944
$ws.Cells.Item(28, 3).Value = 1173.579190908214
$ws.Cells.Item(29, 3).Value = 1164.991612211841
$ws.Cells.Item(30, 3).Value = 1158.999473105604
$ws.Cells.Item(31, 3).Value = 1148.734789083394
$ws.Cells.Item(32, 3).Value = 1142.258602354507
$ws.Cells.Item(33, 3).Value = 1136.165697202429
$ws.Cells.Item(34, 3).Value = 1130.58398511956
$ws.Cells.Item(35, 3).Value = 1126.168521739079
$ws.Cells.Item(36, 3).Value = 1122.035801938925
$ws.Cells.Item(37, 3).Value = 1117.37621194663
$ws.Cells.Item(38, 3).Value = 1114.593978850609
$ws.Cells.Item(39, 3).Value = 1110.691980493739
$ws.Cells.Item(40, 3).Value = 1105.135606479957
$ws.Cells.Item(41, 3).Value = 1100.800437609961
$ws.Cells.Item(42, 3).Value = 1099.577742572855
$ws.Cells.Item(43, 3).Value = 1096.095585613077
$ws.Cells.Item(44, 3).Value = 1093.627903579919
$ws.Cells.Item(45, 3).Value = 1091.468802174811
$ws.Cells.Item(46, 3).Value = 1090.220435641887
$ws.Cells.Item(47, 3).Value = 1088.356404980839
$ws.Cells.Item(48, 3).Value = 1085.855816165715
$ws.Cells.Item(49, 3).Value = 1084.603518512472
$ws.Cells.Item(50, 3).Value = 1082.109415063703
$ws.Cells.Item(51, 3).Value = 1079.034156118763
$ws.Cells.Item(52, 3).Value = 1075.800249509706
$ws.Cells.Item(53, 3).Value = 1070.235660835417
$ws.Cells.Item(54, 3).Value = 1065.252175556308
$ws.Cells.Item(55, 3).Value = 1056.981261768786
$ws.Cells.Item(56, 3).Value = 1048.976162036141
$ws.Cells.Item(57, 3).Value = 1045.131567055102
$ws.Cells.Item(58, 3).Value = 1042.933755553088
$ws.Cells.Item(59, 3).Value = 1040.393183872723
$ws.Cells.Item(60, 3).Value = 1029.997994659415
$ws.Cells.Item(61, 3).Value = 1010.707575921699
$ws.Cells.Item(62, 3).Value = 1008.525403170424
$ws.Cells.Item(63, 3).Value = 1007.616789987065
$ws.Cells.Item(64, 3).Value = 1005.919432667368
$ws.Cells.Item(65, 3).Value = 1004.834113027875

$wb.Save()
